$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2388.889
$ws.Range("I40").Value = 2333.3333
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 2333.3333
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -2158.3333
$ws.Range("N40").Value = -2850

$ws.Range("H41").Value = 599.5
$ws.Range("I41").Value = 666
$ws.Range("J41").Value = 400
$ws.Range("K41").Value = 666
$ws.Range("L41").Value = 400
$ws.Range("M41").Value = -226
$ws.Range("N41").Value = -1280

$ws.Range("H43").Value = 8250
$ws.Range("I43").Value = 8001
$ws.Range("J43").Value = 8499
$ws.Range("K43").Value = 8001
$ws.Range("L43").Value = 8499
$ws.Range("M43").Value = -7932
$ws.Range("N43").Value = -8637

$ws.Range("H80").Value = 270.77274
$ws.Range("I80").Value = 305.4
$ws.Range("J80").Value = 241.91667
$ws.Range("K80").Value = 916.1999999999999
$ws.Range("L80").Value = 725.75001
$ws.Range("M80").Value = 81.80000000000007
$ws.Range("N80").Value = -2721.75001

$ws.Range("H83").Value = 270.77274
$ws.Range("I83").Value = 305.4
$ws.Range("J83").Value = 241.91667
$ws.Range("K83").Value = 2748.6
$ws.Range("L83").Value = 2177.25003
$ws.Range("M83").Value = 2243.4
$ws.Range("N83").Value = -12161.25003

$ws.Range("H98").Value = 928.1667
$ws.Range("I98").Value = 550.3333
$ws.Range("J98").Value = 2061.6667
$ws.Range("K98").Value = 550.3333
$ws.Range("L98").Value = 2061.6667
$ws.Range("M98").Value = 947.6667
$ws.Range("N98").Value = -5057.6667

$ws.Range("H122").Value = 928.1667
$ws.Range("I122").Value = 550.3333
$ws.Range("J122").Value = 2061.6667
$ws.Range("K122").Value = 1650.9999
$ws.Range("L122").Value = 6185.000100000001
$ws.Range("M122").Value = 799.0001
$ws.Range("N122").Value = -11085.0001

$ws.Range("H131").Value = 3863.3333
$ws.Range("I131").Value = 3545
$ws.Range("J131").Value = 4500
$ws.Range("K131").Value = 10635
$ws.Range("L131").Value = 13500
$ws.Range("M131").Value = -5595
$ws.Range("N131").Value = -23580

$ws.Range("H136").Value = 89999.5
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 89999.5
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 89999.5
$ws.Range("N136").Value = -100199.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2049.5715
$ws.Range("I2").Value = 1822.6154
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 1822.6154
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -1709.6154
$ws.Range("N2").Value = -5226

$ws.Range("H32").Value = 8516.974
$ws.Range("I32").Value = 6020.1816
$ws.Range("J32").Value = 24995.8
$ws.Range("K32").Value = 6020.1816
$ws.Range("L32").Value = 24995.8
$ws.Range("M32").Value = -5733.1816
$ws.Range("N32").Value = -25569.8

$ws.Range("H61").Value = 3106.5789
$ws.Range("I61").Value = 3000.8235
$ws.Range("J61").Value = 4005.5
$ws.Range("K61").Value = 3000.8235
$ws.Range("L61").Value = 4005.5
$ws.Range("M61").Value = -2788.8235
$ws.Range("N61").Value = -4429.5

$ws.Range("H102").Value = 1371.25
$ws.Range("I102").Value = 1336.8182
$ws.Range("J102").Value = 1750
$ws.Range("K102").Value = 1336.8182
$ws.Range("L102").Value = 1750
$ws.Range("M102").Value = 285.1818000000001
$ws.Range("N102").Value = -4994

$ws.Range("H116").Value = 2049.5715
$ws.Range("I116").Value = 1822.6154
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 1822.6154
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = 471.3846000000001
$ws.Range("N116").Value = -9588

$ws.Range("H136").Value = 3106.5789
$ws.Range("I136").Value = 3000.8235
$ws.Range("J136").Value = 4005.5
$ws.Range("K136").Value = 9002.470499999999
$ws.Range("L136").Value = 12016.5
$ws.Range("M136").Value = -6452.470499999999
$ws.Range("N136").Value = -17116.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2049.5715
$ws.Range("I3").Value = 1822.6154
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 1822.6154
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -1708.6154
$ws.Range("N3").Value = -5228

$ws.Range("H86").Value = 4891.5713
$ws.Range("I86").Value = 4358.75
$ws.Range("J86").Value = 5602
$ws.Range("K86").Value = 4358.75
$ws.Range("L86").Value = 5602
$ws.Range("M86").Value = -3235.75
$ws.Range("N86").Value = -7848

$ws.Range("H89").Value = 4891.5713
$ws.Range("I89").Value = 4358.75
$ws.Range("J89").Value = 5602
$ws.Range("K89").Value = 21793.75
$ws.Range("L89").Value = 28010
$ws.Range("M89").Value = -16177.75
$ws.Range("N89").Value = -39242

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 500
$ws.Range("I16").Value = 500
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 500
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -213

$ws.Range("H31").Value = 4027.95
$ws.Range("I31").Value = 1451.875
$ws.Range("J31").Value = 5745.3335
$ws.Range("K31").Value = 1451.875
$ws.Range("L31").Value = 5745.3335
$ws.Range("M31").Value = -1156.875
$ws.Range("N31").Value = -6335.3335

$ws.Range("H34").Value = 4027.95
$ws.Range("I34").Value = 1451.875
$ws.Range("J34").Value = 5745.3335
$ws.Range("K34").Value = 1451.875
$ws.Range("L34").Value = 5745.3335
$ws.Range("M34").Value = -1249.875
$ws.Range("N34").Value = -6149.3335

$ws.Range("H62").Value = 34513.23
$ws.Range("I62").Value = 4056.0833
$ws.Range("J62").Value = 399999
$ws.Range("K62").Value = 4056.0833
$ws.Range("L62").Value = 399999
$ws.Range("M62").Value = -3432.0833
$ws.Range("N62").Value = -401247

$ws.Range("H65").Value = 34513.23
$ws.Range("I65").Value = 4056.0833
$ws.Range("J65").Value = 399999
$ws.Range("K65").Value = 20280.4165
$ws.Range("L65").Value = 1999995
$ws.Range("M65").Value = -17160.4165
$ws.Range("N65").Value = -2006235

$ws.Range("H99").Value = 14405.619
$ws.Range("I99").Value = 11462.556
$ws.Range("J99").Value = 16612.916
$ws.Range("K99").Value = 11462.556
$ws.Range("L99").Value = 16612.916
$ws.Range("M99").Value = -9964.556
$ws.Range("N99").Value = -19608.916

$ws.Range("H113").Value = 500
$ws.Range("I113").Value = 500
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 500
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1670

$ws.Range("H126").Value = 14405.619
$ws.Range("I126").Value = 11462.556
$ws.Range("J126").Value = 16612.916
$ws.Range("K126").Value = 34387.66800000001
$ws.Range("L126").Value = 49838.74800000001
$ws.Range("M126").Value = -31917.66800000001
$ws.Range("N126").Value = -54778.74800000001

$ws.Range("H132").Value = 3698.25
$ws.Range("I132").Value = 2847.6667
$ws.Range("J132").Value = 6250
$ws.Range("K132").Value = 8543.000100000001
$ws.Range("L132").Value = 18750
$ws.Range("M132").Value = -6013.000100000001
$ws.Range("N132").Value = -23810

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 300
$ws.Range("I132").Value = 300
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2700
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -170

$ws.Range("H134").Value = 3004.5
$ws.Range("I134").Value = 3004.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9013.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3943.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 12188.4
$ws.Range("I99").Value = 2147.3333
$ws.Range("J99").Value = 27250
$ws.Range("K99").Value = 2147.3333
$ws.Range("L99").Value = 27250
$ws.Range("M99").Value = 98.66670000000022
$ws.Range("N99").Value = -31742

$ws.Range("H102").Value = 1928.9714
$ws.Range("I102").Value = 1321.9474
$ws.Range("J102").Value = 2649.8125
$ws.Range("K102").Value = 1321.9474
$ws.Range("L102").Value = 2649.8125
$ws.Range("M102").Value = 300.0526
$ws.Range("N102").Value = -5893.8125

$ws.Range("H132").Value = 2368.743
$ws.Range("I132").Value = 1608.3889
$ws.Range("J132").Value = 3173.8235
$ws.Range("K132").Value = 4825.1667
$ws.Range("L132").Value = 9521.470499999999
$ws.Range("M132").Value = -2295.1667
$ws.Range("N132").Value = -14581.4705

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3273.2666
$ws.Range("I46").Value = 2485.4285
$ws.Range("J46").Value = 3962.625
$ws.Range("K46").Value = 2485.4285
$ws.Range("L46").Value = 3962.625
$ws.Range("M46").Value = -2297.4285
$ws.Range("N46").Value = -4338.625

$ws.Range("H55").Value = 412.78946
$ws.Range("I55").Value = 365.1875
$ws.Range("J55").Value = 666.6667
$ws.Range("K55").Value = 365.1875
$ws.Range("L55").Value = 666.6667
$ws.Range("M55").Value = -192.1875
$ws.Range("N55").Value = -1012.6667

$ws.Range("H61").Value = 4298.636
$ws.Range("I61").Value = 4298.636
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4298.636
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4096.636

$ws.Range("H113").Value = 4298.636
$ws.Range("I113").Value = 4298.636
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4298.636
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2128.636

$ws.Range("H132").Value = 5960
$ws.Range("I132").Value = 5000
$ws.Range("J132").Value = 6200
$ws.Range("K132").Value = 15000
$ws.Range("L132").Value = 18600
$ws.Range("M132").Value = -12470
$ws.Range("N132").Value = -23660

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3247.6667
$ws.Range("I122").Value = 3435.8462
$ws.Range("J122").Value = 2024.5
$ws.Range("K122").Value = 10307.5386
$ws.Range("L122").Value = 6073.5
$ws.Range("M122").Value = -7857.5386
$ws.Range("N122").Value = -10973.5

$ws.Range("H132").Value = 1835
$ws.Range("I132").Value = 1750
$ws.Range("J132").Value = 2005
$ws.Range("K132").Value = 5250
$ws.Range("L132").Value = 6015
$ws.Range("M132").Value = -2720
$ws.Range("N132").Value = -11075

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N132").ClearContents()
